$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3
$ws.Range("H6").Value = 4
$ws.Range("H10").Value = 3
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 3
$ws.Range("H15").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 2
$ws.Range("H22").Value = 2
$ws.Range("H23").Value = 2
$ws.Range("H24").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1

$ws.Range("H28").Select()
